$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5
$ws.Range("P2").Value = 0.5

$ws.Range("P4").Value = 1

$ws.Range("J6").Value = 0.3
$ws.Range("O6").Value = 0.1
$ws.Range("Q6").Value = 0.2
$ws.Range("S6").Value = 0.4

$ws.Range("F7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.2222222222222222
$ws.Range("S7").Value = 0.6666666666666666

$ws.Range("B8").Value = 0.04761904761904762
$ws.Range("D8").Value = 0.04761904761904762
$ws.Range("F8").Value = 0.04761904761904762
$ws.Range("J8").Value = 0.04761904761904762
$ws.Range("O8").Value = 0.04761904761904762
$ws.Range("Q8").Value = 0.04761904761904762
$ws.Range("R8").Value = 0.1428571428571428
$ws.Range("S8").Value = 0.5714285714285714

$ws.Range("D9").Value = 0.08333333333333333
$ws.Range("J9").Value = 0.1666666666666667
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.5

$ws.Range("B10").Value = 0.06451612903225806
$ws.Range("D10").Value = 0.03225806451612903
$ws.Range("F10").Value = 0.1290322580645161
$ws.Range("J10").Value = 0.03225806451612903
$ws.Range("Q10").Value = 0.1612903225806452
$ws.Range("R10").Value = 0.06451612903225806
$ws.Range("S10").Value = 0.5161290322580645

$ws.Range("G11").Value = 0.1111111111111111
$ws.Range("K11").Value = 0.1111111111111111
$ws.Range("L11").Value = 0.7777777777777778

$ws.Range("G12").Value = 1

$ws.Range("G13").Value = 1

$ws.Range("I15").Value = 0.25
$ws.Range("J15").Value = 0.25
$ws.Range("O15").Value = 0.125
$ws.Range("S15").Value = 0.375

$ws.Range("H16").Value = 0.4285714285714285
$ws.Range("J16").Value = 0.2857142857142857
$ws.Range("M16").Value = 0.1428571428571428
$ws.Range("O16").Value = 0.1428571428571428

$ws.Range("F17").Value = 0.2
$ws.Range("H17").Value = 0.1
$ws.Range("I17").Value = 0.3
$ws.Range("J17").Value = 0.2
$ws.Range("K17").Value = 0.2

$ws.Range("F18").Value = 0.1666666666666667
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.1666666666666667
$ws.Range("S18").Value = 0.1666666666666667

$ws.Range("H19").Value = 0.3529411764705883
$ws.Range("I19").Value = 0.1372549019607843
$ws.Range("J19").Value = 0.2549019607843137
$ws.Range("K19").Value = 0.09803921568627451
$ws.Range("M19").Value = 0.0196078431372549
$ws.Range("O19").Value = 0.05882352941176471
$ws.Range("S19").Value = 0.07843137254901961
